$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: Replace the 5 "características a testar" bullet paragraphs with
# 11 new bullet paragraphs.
# ---------------------------------------------------------------------------

$newTexts1 = @(
    "*Teste de integração e de segurança do banco de dados. (M)",
    "*Teste unitário do backend do controle de estoque. (A)",
    "*Teste unitário do backend do controle de cadastro. (A)",
    "*Teste unitário do backend do controle de orçamento. (A)",
    "*Teste unitário do backend do controle de agendamento. (A)",
    "*Teste unitário do frontend do controle de estoque. (A)",
    "*Teste unitário do frontend do controle de cadastro. (A)",
    "*Teste unitário do frontend do controle de orçamento. (A)",
    "*Teste unitário do frontend do controle de agendamento. (A)",
    "*Teste de segurança do login do website. (M)",
    "*Teste de integração do host do website. (A)"
)

$startIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.StartsWith("1. Funcionalidade do banco de dados")) {
        $startIndex = $i
        break
    }
}
Write-Host "Start index: " $startIndex

# Replace the text of the existing 5 paragraphs with the first 5 new texts.
for ($k = 0; $k -lt 5; $k++) {
    $p = $d.Paragraphs.Item($startIndex + $k)
    $p.Range.Text = $newTexts1[$k]
}

# Insert 6 additional empty paragraphs right after the 5th (now last) paragraph
# of this block, then fill them in with the remaining new texts.
$lastOfBlock = $d.Paragraphs.Item($startIndex + 4)
$insertRange = $lastOfBlock.Range
for ($k = 0; $k -lt 6; $k++) {
    $insertRange.InsertParagraphAfter()
}
for ($k = 0; $k -lt 6; $k++) {
    $p = $d.Paragraphs.Item($startIndex + 5 + $k)
    $p.Range.Text = $newTexts1[5 + $k]
}

# ---------------------------------------------------------------------------
# Part 2: Collapse the "Os seguintes testes..." intro paragraph plus the
# following 7 bullet paragraphs into a single explanatory paragraph.
# ---------------------------------------------------------------------------

$mergedText = "Todos os componentes desenvolvidos pela equipe precisam ser testados de uma forma ou outra, com a exceção do sistema de notas fiscais, pois a equipe ainda não tem certeza se esse componente será desenvolvido, caso contrário o seu teste será adicionado na sprint."

$introIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.StartsWith("Os seguintes testes deverão ser realizados")) {
        $introIndex = $i
        break
    }
}
Write-Host "Intro index: " $introIndex

$pIntro = $d.Paragraphs.Item($introIndex)
$pIntro.Range.Text = $mergedText

# Insert a fresh empty paragraph right after the last of the 7 bullet
# paragraphs that follow (this becomes the carrier for the trailing empty
# run left behind once those bullet paragraphs are deleted and it gets
# merged back up into the text paragraph).
$lastBullet = $d.Paragraphs.Item($introIndex + 7)
$lastBullet.Range.InsertParagraphAfter()

# Remove the 7 old bullet paragraphs that now sit between the text
# paragraph and the new carrier paragraph.
for ($k = 0; $k -lt 7; $k++) {
    $d.Paragraphs.Item($introIndex + 1).Range.Delete()
}

# Merge the now-empty trailing paragraph back into the text paragraph by
# removing the paragraph mark that separates them - this leaves the empty
# run behind inside the text paragraph, matching the target structure.
$pIntro3 = $d.Paragraphs.Item($introIndex)
$markPos = $pIntro3.Range.End - 1
$markRange = $d.Range($markPos, $markPos + 1)
$markRange.Delete()

Write-Host "Done."
